$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.69462833333333
$ws.Range("H2").Value = 164.083885
$ws.Range("I2").Value = 0.2790924419198448
$ws.Range("J2").Value = 0.2790924419198448
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 482.4420275560744
$ws.Range("R2").Value = 4341.978248004671
$ws.Range("S2").Value = 0.01790418946728285
$ws.Range("T2").Value = 0.01790418946728285
$ws.Range("G3").Value = 54.69462833333333
$ws.Range("H3").Value = 164.083885
$ws.Range("I3").Value = 0.2790924419198448
$ws.Range("J3").Value = 0.2790924419198448
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 2992.427072371158
$ws.Range("R3").Value = 26931.84365134043
$ws.Range("S3").Value = 0.111053718810873
$ws.Range("T3").Value = 0.111053718810873
$ws.Range("G4").Value = 54.69462833333333
$ws.Range("H4").Value = 164.083885
$ws.Range("I4").Value = 0.2790924419198448
$ws.Range("J4").Value = 0.2790924419198448
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 1198.258850982628
$ws.Range("R4").Value = 10784.32965884365
$ws.Range("S4").Value = 0.04446928806663308
$ws.Range("T4").Value = 0.04446928806663308
$ws.Range("G5").Value = 54.69462833333333
$ws.Range("H5").Value = 164.083885
$ws.Range("I5").Value = 0.2790924419198448
$ws.Range("J5").Value = 0.2790924419198448
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 2847.230555205738
$ws.Range("R5").Value = 25625.07499685164
$ws.Range("S5").Value = 0.1056652455750558
$ws.Range("T5").Value = 0.1056652455750559
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09859081282432611
$ws.Range("J6").Value = 0.09859081282432611
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 170.4250796265912
$ws.Range("R6").Value = 1533.82571663932
$ws.Range("S6").Value = 0.006324745236372669
$ws.Range("T6").Value = 0.006324745236372669
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09859081282432611
$ws.Range("J7").Value = 0.09859081282432611
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("R7").Value = 9513.809611533301
$ws.Range("S7").Value = 0.03923028631449874
$ws.Range("T7").Value = 0.03923028631449875
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09859081282432611
$ws.Range("J8").Value = 0.09859081282432611
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 423.2909830150446
$ws.Range("R8").Value = 3809.618847135401
$ws.Range("S8").Value = 0.01570900030846274
$ws.Range("T8").Value = 0.01570900030846274
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09859081282432611
$ws.Range("J9").Value = 0.09859081282432611
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 1005.798554790714
$ws.Range("R9").Value = 9052.18699311642
$ws.Range("S9").Value = 0.03732678096499194
$ws.Range("T9").Value = 0.03732678096499195
$ws.Range("G10").Value = 11.023718
$ws.Range("H10").Value = 33.071154
$ws.Range("I10").Value = 0.05625116157486912
$ws.Range("J10").Value = 0.05625116157486911
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 97.23632878011867
$ws.Range("R10").Value = 875.126959021068
$ws.Range("S10").Value = 0.003608594513213098
$ws.Range("T10").Value = 0.003608594513213098
$ws.Range("G11").Value = 11.023718
$ws.Range("H11").Value = 33.071154
$ws.Range("I11").Value = 0.05625116157486912
$ws.Range("J11").Value = 0.05625116157486911
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 603.1245331871301
$ws.Range("R11").Value = 5428.12079868417
$ws.Range("S11").Value = 0.02238290882171078
$ws.Range("T11").Value = 0.02238290882171078
$ws.Range("G12").Value = 11.023718
$ws.Range("H12").Value = 33.071154
$ws.Range("I12").Value = 0.05625116157486912
$ws.Range("J12").Value = 0.05625116157486911
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 241.5094144846067
$ws.Range("R12").Value = 2173.58473036146
$ws.Range("S12").Value = 0.008962797741667227
$ws.Range("T12").Value = 0.008962797741667227
$ws.Range("G13").Value = 11.023718
$ws.Range("H13").Value = 33.071154
$ws.Range("I13").Value = 0.05625116157486912
$ws.Range("J13").Value = 0.05625116157486911
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 573.8601335817619
$ws.Range("R13").Value = 5164.741202235858
$ws.Range("S13").Value = 0.021296860498278
$ws.Range("T13").Value = 0.021296860498278
$ws.Range("G14").Value = 110.9336623333333
$ws.Range("H14").Value = 332.800987
$ws.Range("I14").Value = 0.5660655836809599
$ws.Range("J14").Value = 0.5660655836809599
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 978.506712837417
$ws.Range("R14").Value = 8806.560415536753
$ws.Range("S14").Value = 0.03631393738725003
$ws.Range("T14").Value = 0.03631393738725003
$ws.Range("G15").Value = 110.9336623333333
$ws.Range("H15").Value = 332.800987
$ws.Range("I15").Value = 0.5660655836809599
$ws.Range("J15").Value = 0.5660655836809599
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 6069.350949428347
$ws.Range("R15").Value = 54624.15854485513
$ws.Range("S15").Value = 0.2252432481732072
$ws.Range("T15").Value = 0.2252432481732072
$ws.Range("G16").Value = 110.9336623333333
$ws.Range("H16").Value = 332.800987
$ws.Range("I16").Value = 0.5660655836809599
$ws.Range("J16").Value = 0.5660655836809599
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 2430.352793563514
$ws.Range("R16").Value = 21873.17514207163
$ws.Range("S16").Value = 0.0901942500920356
$ws.Range("T16").Value = 0.0901942500920356
$ws.Range("G17").Value = 110.9336623333333
$ws.Range("H17").Value = 332.800987
$ws.Range("I17").Value = 0.5660655836809599
$ws.Range("J17").Value = 0.5660655836809599
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 5774.858018439943
$ws.Range("R17").Value = 51973.72216595949
$ws.Range("S17").Value = 0.214314148028467
$ws.Range("T17").Value = 0.214314148028467